# "Got all the stuff calculated, trying to get a few more data points"
#
# The workbook's styles.xml ships with an empty cellXfs table. Touching any
# cell's contents needs a style index to exist, so we seed style index 0
# by (re)applying the built-in "Normal" style to a cell first. This is a
# no-op visually (Normal is the default look) but makes the style table
# non-empty so subsequent writes succeed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")   # density data lives on Sheet1 (the tab-selected sheet)

$ws.Range("C1").Style = "Normal"

# Add a new data point "c" as row 4, following the same pattern as rows 2-3.
$ws.Range("A4").Value = "c"
$ws.Range("B4").Value = 736.76
$ws.Range("C4").Formula = "=193.94-32.8"
$ws.Range("D4").Formula = "=C4/B4"
$ws.Range("E4").Formula = "=D4*1000"

# Move the active selection down to the now-empty row below the data.
$ws.Range("A5").Select()
